$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (moon / first-animal question) - update backward reasoning columns
$ws.Range("F2").Value = "The conclusion is plausible because humans are technically animals, and thus, Neil Armstrong, as a human, could be considered the first animal sent to the Moon based on the premise."
$ws.Range("G2").Value = "The conclusion that humans, being animals, could be considered the first animal sent to the Moon, is plausible. This is because the premise states that no animals were ever sent to the Moon, and it is a fact that humans belong to the species Homo sapiens. Therefore, the conclusion logically follows from the given premise and aligns with our broader understanding of the classification of humans within the animal kingdom."
$ws.Range("H2").Value = "Who was technically the first animal sent to the Moon according to the given information?"

# Row 3 (Leonardo DiCaprio question) - update forward answer + backward reasoning columns + original response
$ws.Range("D3").Value = "Leonardo DiCaprio does not have any children."
$ws.Range("F3").Value = "It is highly plausible that Leonardo DiCaprio does not have any children."
$ws.Range("G3").Value = "Leonardo DiCaprio has never been married, and there is no publicly available information or evidence to suggest that he has children. The lack of any reports, public acknowledgment, or visible presence of children strongly suggests that the conclusion about him not having any children is credible. This aligns with his public persona and is supported both by specific pieces of evidence and general knowledge about his personal life."
$ws.Range("H3").Value = "Does Leonardo DiCaprio have any children?"
$ws.Range("I3").Value = "{'Answer:': 'Leonardo DiCaprio does not have any children.', 'Source:': 'Quora (https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have)', 'Premise of the Question:': 'InvalidLeonardo DiCaprio does not have any children, so the question contains a false premise.', 'Explanation:': ''}"

$wb.Save()
